$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# --- Fix typo in existing command description (Z4, "web" list) ---
$ws.Range("Z4").Value = 'assertAttributeContain(locator,attrName,contains)'

# --- Update the "io" command list (column L) to add assertPath(path) and
#     rename saveMatches(...) to its new signature, shifting the rest down ---
$ws.Range("L4").Value = 'assertPath(path)'
$ws.Range("L5").Value = 'assertReadableFile(file,minByte)'
$ws.Range("L6").Value = 'base64(var,file)'
$ws.Range("L7").Value = 'compare(expected,actual,failFast)'
$ws.Range("L8").Value = 'copyFiles(source,target)'
$ws.Range("L9").Value = 'copyFilesByRegex(sourceDir,regex,target)'
$ws.Range("L10").Value = 'count(var,path,pattern)'
$ws.Range("L11").Value = 'deleteFiles(location,recursive)'
$ws.Range("L12").Value = 'deleteFilesByRegex(sourceDir,regex)'
$ws.Range("L13").Value = 'filter(source,target,matchPattern)'
$ws.Range("L14").Value = 'makeDirectory(source)'
$ws.Range("L15").Value = 'moveFiles(source,target)'
$ws.Range("L16").Value = 'moveFilesByRegex(sourceDir,regex,target)'
$ws.Range("L17").Value = 'readFile(var,file)'
$ws.Range("L18").Value = 'readProperty(var,file,property)'
$ws.Range("L19").Value = 'rename(target,newName)'
$ws.Range("L20").Value = 'saveDiff(var,expected,actual)'
$ws.Range("L21").Value = 'saveFileMeta(var,file)'
$ws.Range("L22").Value = 'saveMatches(var,path,fileFilter,textFilter)'
$ws.Range("L23").Value = 'searchAndReplace(file,config,saveAs)'
$ws.Range("L24").Value = 'unzip(zipFile,target)'
$ws.Range("L25").Value = 'validate(var,profile,inputFile)'
$ws.Range("L26").Value = 'writeBase64decode(encodedSource,decodedTarget,append)'
$ws.Range("L27").Value = 'writeFile(file,content,append)'
$ws.Range("L28").Value = 'writeFileAsIs(file,content,append)'
$ws.Range("L29").Value = 'writeProperty(file,property,value)'
$ws.Range("L30").Value = 'zip(filePattern,zipFile)'
$ws.Range("F19").Value = 'clearVariables(variables)'
$ws.Range("F20").Value = 'failImmediate(text)'
$ws.Range("F21").Value = 'incrementChar(var,amount,config)'
$ws.Range("F22").Value = 'macro(file,sheet,name)'
$ws.Range("F23").Value = 'outputToCloud(resource)'
$ws.Range("F24").Value = 'prependText(var,prependWith)'
$ws.Range("F25").Value = 'repeatUntil(steps,maxWaitMs)'
$ws.Range("F26").Value = 'save(var,value)'
$ws.Range("F27").Value = 'saveCount(text,regex,saveVar)'
$ws.Range("F28").Value = 'saveMatches(text,regex,saveVar)'
$ws.Range("F29").Value = 'saveReplace(text,regex,replace,saveVar)'
$ws.Range("F30").Value = 'saveVariablesByPrefix(var,prefix)'
$ws.Range("F31").Value = 'saveVariablesByRegex(var,regex)'
$ws.Range("F32").Value = 'section(steps)'
$ws.Range("F33").Value = 'split(text,delim,saveVar)'
$ws.Range("F34").Value = 'startRecording()'
$ws.Range("F35").Value = 'stopRecording()'
$ws.Range("F36").Value = 'substringAfter(text,delim,saveVar)'
$ws.Range("F37").Value = 'substringBefore(text,delim,saveVar)'
$ws.Range("F38").Value = 'substringBetween(text,start,end,saveVar)'
$ws.Range("F39").Value = 'verbose(text)'
$ws.Range("F40").Value = 'waitFor(waitMs)'

# --- Update named ranges to reflect the new list sizes ---
$wb.Names.Item("base").RefersTo = '=''#system''!$F$2:$F$40'
$wb.Names.Item("io").RefersTo = '=''#system''!$L$2:$L$30'
